# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" colours (bound to the Notes Master)
#   ppt/theme/theme2.xml  -> "Integral" colours      (bound to the Slide Master
#                                                      that every slide/layout uses)
# The authored edit swaps the two themes' colour content: the theme that
# drives the slide master/slides switches from the "Integral" palette to the
# stock "Office Theme" palette.
#
# PowerPoint's Design > Variants > Colors > Customize Colors dialog edits the
# live theme colour scheme in place (12 slots, in clrScheme document order:
# Dark1, Light1, Dark2, Light2, Accent1-6, Hyperlink, FollowedHyperlink).
# That is exactly the ThemeColorScheme collection on the active theme.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

function ToRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the stock Office "Office Theme" colours.
$colors.Item(1).RGB  = ToRGB 0x00 0x00 0x00   # Dark 1    (dk1)      000000
$colors.Item(2).RGB  = ToRGB 0xFF 0xFF 0xFF   # Light 1   (lt1)      FFFFFF
$colors.Item(3).RGB  = ToRGB 0x44 0x54 0x6A   # Dark 2    (dk2)      44546A
$colors.Item(4).RGB  = ToRGB 0xE7 0xE6 0xE6   # Light 2   (lt2)      E7E6E6
$colors.Item(5).RGB  = ToRGB 0x5B 0x9B 0xD5   # Accent 1             5B9BD5
$colors.Item(6).RGB  = ToRGB 0xED 0x7D 0x31   # Accent 2             ED7D31
$colors.Item(7).RGB  = ToRGB 0xA5 0xA5 0xA5   # Accent 3             A5A5A5
$colors.Item(8).RGB  = ToRGB 0xFF 0xC0 0x00   # Accent 4             FFC000
$colors.Item(9).RGB  = ToRGB 0x44 0x72 0xC4   # Accent 5             4472C4
$colors.Item(10).RGB = ToRGB 0x70 0xAD 0x47   # Accent 6             70AD47
$colors.Item(11).RGB = ToRGB 0x05 0x63 0xC1   # Hyperlink            0563C1
$colors.Item(12).RGB = ToRGB 0x95 0x4F 0x72   # Followed Hyperlink   954F72
